{"js": "// Replace the invoice placeholder values in the table cells.\n// Each pair is [old exact text, new exact text]; all are unique substrings\n// within the document body (the two \"wq qw wq qw\" cells are left untouched).\nconst replacements = [\n  [\"wq 23/32q\", \"dede dedew/2e\"],\n  [\"1111111111 eqe\", \"323232 eqe\"],\n  [\"wdw\", \"dede\"],\n  [\"dwd@dede\", \"de@ded\"],\n  [\"+380984343994\", \"+380987676554\"],\n  [\"dedeed\", \"ddcdc\"],\n  [\"Amount USD: -610080234\", \"Amount USD: 254\"],\n  [\"___________________16-4-2020\", \"___________________12-4-2020\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the invoice placeholder values in the table cells.\n# Each pair is (old exact text, new exact text); all are unique substrings\n# within the document body (the two \"wq qw wq qw\" cells are left untouched).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"wq 23/32q\", \"dede dedew/2e\"),\n    @(\"1111111111 eqe\", \"323232 eqe\"),\n    @(\"wdw\", \"dede\"),\n    @(\"dwd@dede\", \"de@ded\"),\n    @(\"+380984343994\", \"+380987676554\"),\n    @(\"dedeed\", \"ddcdc\"),\n    @(\"Amount USD: -610080234\", \"Amount USD: 254\"),\n    @(\"___________________16-4-2020\", \"___________________12-4-2020\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue = 1 (Wrap), wdReplaceAll = 2 (Replace)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
